$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "2023-08-25 07:27:12"
$ws.Range("B9").Value = 6
